$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.352.18'
$ws.Range("E2").Value = '  -0.32%  '

$ws.Range("D3").Value = '2.931.24'
$ws.Range("E3").Value = '  +0.90%  '

$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.36%  '

$ws.Range("E5").Value = '  +5.95%  '

$ws.Range("D6").Value = '104.48'
$ws.Range("E6").Value = '  -1.67%  '

$ws.Range("E7").Value = '  -1.54%  '

$ws.Range("D8").Value = '''1.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.12%  '

$ws.Range("E9").Value = '  -2.71%  '

$ws.Range("D10").Value = '36.84'
$ws.Range("E10").Value = '  -1.62%  '

$ws.Range("D11").Value = '''0.140'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.40%  '

$ws.Range("D12").Value = '0.0838'
$ws.Range("E12").Value = '  -0.50%  '

$ws.Range("D13").Value = '18.42'
$ws.Range("E13").Value = '  -2.20%  '

$ws.Range("D14").Value = '3.391.57'
$ws.Range("E14").Value = '  +0.49%  '

$ws.Range("D15").Value = '7.41'
$ws.Range("E15").Value = '  -1.26%  '

$ws.Range("D16").Value = '2.919.78'
$ws.Range("E16").Value = '  -0.06%  '

$ws.Range("D17").Value = '0.943'
$ws.Range("E17").Value = '  -1.59%  '

$ws.Range("D18").Value = '51.257.30'
$ws.Range("E18").Value = '  -0.50%  '

$ws.Range("D19").Value = '3.26'
$ws.Range("E19").Value = '  -4.98%  '

$ws.Range("E20").Value = '  -0.67%  '

$ws.Range("D21").Value = '12.97'
$ws.Range("E21").Value = '  -2.48%  '

$ws.Range("D22").Value = '0.0₃0945'
$ws.Range("E22").Value = '  -1.18%  '

$ws.Range("D23").Value = '''68.40'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.30%  '

$ws.Range("D24").Value = '260.59'
$ws.Range("E24").Value = '  +0.14%  '

$ws.Range("E25").Value = '  +0.50%  '

$ws.Range("D26").Value = '4.35'
$ws.Range("E26").Value = '  +4.40%  '

$ws.Range("D27").Value = '0.174'
$ws.Range("E27").Value = '  +2.59%  '

$ws.Range("D28").Value = '''1.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.15%  '

$ws.Range("E29").Value = '  -1.87%  '

$ws.Range("E30").Value = '  -7.09%  '

$ws.Range("E31").Value = '  +0.72%  '

$ws.Range("E32").Value = '  +3.01%  '

$ws.Range("D33").Value = '9.94'
$ws.Range("E33").Value = '  -2.02%  '

$ws.Range("E34").Value = '  -1.12%  '

$ws.Range("E35").Value = '  -0.93%  '

$ws.Range("D36").Value = '50.93'
$ws.Range("E36").Value = '  +0.57%  '

$ws.Range("D37").Value = '''1.00'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.27%  '

$ws.Range("D38").Value = '0.0424'
$ws.Range("E38").Value = '  +0.68%  '

$ws.Range("D39").Value = '3.04'
$ws.Range("E39").Value = '  -1.51%  '

$ws.Range("E40").Value = '  +1.84%  '

$ws.Range("D41").Value = '17.15'
$ws.Range("E41").Value = '  -1.98%  '

$ws.Range("E42").Value = '  -3.93%  '

$ws.Range("E43").Value = '  -1.50%  '

$ws.Range("D44").Value = '22.31'
$ws.Range("E44").Value = '  +0.48%  '

$ws.Range("D45").Value = '119.29'
$ws.Range("E45").Value = '  +0.45%  '

$ws.Range("E46").Value = '  -2.74%  '

$ws.Range("D47").Value = '2.030.97'
$ws.Range("E47").Value = '  -2.40%  '

$ws.Range("E48").Value = '  +0.49%  '

$ws.Range("D49").Value = '3.18'
$ws.Range("E49").Value = '  -3.17%  '

$ws.Range("E50").Value = '  +4.82%  '

$ws.Range("D51").Value = '3.209.09'
$ws.Range("E51").Value = '  +0.01%  '
